{"js": "// Replace the 25 \"two-digit \u00f7 one-digit\" answer strings in the table\n// with their new values, as described by the diff. Each old value is\n// unique in the document, so a targeted search+replace per pair is safe.\nconst replacements = [\n  [\"11\u00f77=1, 4\", \"72\u00f75=14, 2\"],\n  [\"64\u00f78=8, 0\", \"73\u00f73=24, 1\"],\n  [\"44\u00f75=8, 4\", \"68\u00f77=9, 5\"],\n  [\"77\u00f75=15, 2\", \"81\u00f78=10, 1\"],\n  [\"54\u00f77=7, 5\", \"93\u00f72=46, 1\"],\n  [\"55\u00f73=18, 1\", \"25\u00f78=3, 1\"],\n  [\"58\u00f72=29, 0\", \"23\u00f72=11, 1\"],\n  [\"22\u00f79=2, 4\", \"98\u00f76=16, 2\"],\n  [\"87\u00f72=43, 1\", \"99\u00f73=33, 0\"],\n  [\"41\u00f78=5, 1\", \"25\u00f72=12, 1\"],\n  [\"11\u00f75=2, 1\", \"58\u00f76=9, 4\"],\n  [\"39\u00f73=13, 0\", \"29\u00f79=3, 2\"],\n  [\"71\u00f72=35, 1\", \"92\u00f76=15, 2\"],\n  [\"67\u00f76=11, 1\", \"83\u00f77=11, 6\"],\n  [\"69\u00f73=23, 0\", \"51\u00f78=6, 3\"],\n  [\"19\u00f74=4, 3\", \"42\u00f75=8, 2\"],\n  [\"10\u00f75=2, 0\", \"57\u00f78=7, 1\"],\n  [\"86\u00f73=28, 2\", \"56\u00f72=28, 0\"],\n  [\"40\u00f79=4, 4\", \"31\u00f73=10, 1\"],\n  [\"31\u00f72=15, 1\", \"98\u00f75=19, 3\"],\n  [\"16\u00f78=2, 0\", \"66\u00f73=22, 0\"],\n  [\"22\u00f78=2, 6\", \"16\u00f73=5, 1\"],\n  [\"11\u00f76=1, 5\", \"80\u00f72=40, 0\"],\n  [\"88\u00f78=11, 0\", \"18\u00f79=2, 0\"],\n  [\"90\u00f77=12, 6\", \"58\u00f75=11, 3\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 \"two-digit \u00f7 one-digit\" answer strings in the table\n# with their new values, as described by the diff. Each old value is\n# unique in the document, so Find/Replace (one occurrence each) is safe.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"11\u00f77=1, 4\", \"72\u00f75=14, 2\"),\n    @(\"64\u00f78=8, 0\", \"73\u00f73=24, 1\"),\n    @(\"44\u00f75=8, 4\", \"68\u00f77=9, 5\"),\n    @(\"77\u00f75=15, 2\", \"81\u00f78=10, 1\"),\n    @(\"54\u00f77=7, 5\", \"93\u00f72=46, 1\"),\n    @(\"55\u00f73=18, 1\", \"25\u00f78=3, 1\"),\n    @(\"58\u00f72=29, 0\", \"23\u00f72=11, 1\"),\n    @(\"22\u00f79=2, 4\", \"98\u00f76=16, 2\"),\n    @(\"87\u00f72=43, 1\", \"99\u00f73=33, 0\"),\n    @(\"41\u00f78=5, 1\", \"25\u00f72=12, 1\"),\n    @(\"11\u00f75=2, 1\", \"58\u00f76=9, 4\"),\n    @(\"39\u00f73=13, 0\", \"29\u00f79=3, 2\"),\n    @(\"71\u00f72=35, 1\", \"92\u00f76=15, 2\"),\n    @(\"67\u00f76=11, 1\", \"83\u00f77=11, 6\"),\n    @(\"69\u00f73=23, 0\", \"51\u00f78=6, 3\"),\n    @(\"19\u00f74=4, 3\", \"42\u00f75=8, 2\"),\n    @(\"10\u00f75=2, 0\", \"57\u00f78=7, 1\"),\n    @(\"86\u00f73=28, 2\", \"56\u00f72=28, 0\"),\n    @(\"40\u00f79=4, 4\", \"31\u00f73=10, 1\"),\n    @(\"31\u00f72=15, 1\", \"98\u00f75=19, 3\"),\n    @(\"16\u00f78=2, 0\", \"66\u00f73=22, 0\"),\n    @(\"22\u00f78=2, 6\", \"16\u00f73=5, 1\"),\n    @(\"11\u00f76=1, 5\", \"80\u00f72=40, 0\"),\n    @(\"88\u00f78=11, 0\", \"18\u00f79=2, 0\"),\n    @(\"90\u00f77=12, 6\", \"58\u00f75=11, 3\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
